# Updates the single-column results table in place: several individual
# metric cells get new values, and the last three rows (which still held
# the old full tab-separated summary line) get collapsed down to just
# the first figure from that line.
$d = $word.ActiveDocument
$t = $d.Tables(1)

$changes = @{
    1  = "0M"
    2  = "0M"
    3  = "0M"
    4  = "84"
    6  = "0.00049"
    7  = "0.00016"
    8  = "0.00006"
    9  = "0.00022"
    10 = "0.00025"
    11 = "0.00034"
    12 = "0.01389"
    44 = "99.9"
    45 = "0.01"
    46 = "14"
}

foreach ($rowIndex in $changes.Keys) {
    $cell = $t.Rows($rowIndex).Cells(1)
    $cell.Range.Text = $changes[$rowIndex]
}
